# Design Rationale.docx — apply the "Updated Class Diagrams and Design
# Rationale for Beating up Zombies" edits.
#
# Summary of changes (in document order):
#   1. Bold "DRY" in "...we used the principle DRY (Don't Repeat Yourself)."
#   2. Bold "DRY" in "...using the principle of DRY. Before beginning..."
#   3. "Since classes should be responsible..." -> "Using the design
#      principle that classes should be responsible...", with "classes
#      should be responsible for their own properties" bolded.
#   4. Append a new sentence to the "Beating up the Zombies" paragraph
#      discussing FF (Fail Fast), DRY and grouping/encapsulation, with the
#      relevant principle names/phrases bolded.
#   5. Wrap the "Crafting class has 2 attributes ... Don't Repeat
#      Yourself. " sentence in a "_GoBack" bookmark (last-edit marker).

$d = $word.ActiveDocument

function Find-Unique([string]$searchText) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "anchor text not found: $searchText"
    }
    return $r
}

# ---------------------------------------------------------------------
# Edit 1: bold the first "DRY" — "...we used the principle DRY (Don't..."
# ---------------------------------------------------------------------
$r1 = Find-Unique "principle DRY ("
$dry1Start = $r1.Start + 10
$dry1End = $dry1Start + 3
$d.Range($dry1Start, $dry1End).Bold = 1

# ---------------------------------------------------------------------
# Edit 2: bold the second "DRY" — "...principle of DRY. Before beginning..."
# ---------------------------------------------------------------------
$r2 = Find-Unique "principle of DRY. Before"
$dry2Start = $r2.Start + 13
$dry2End = $dry2Start + 3
$d.Range($dry2Start, $dry2End).Bold = 1

# ---------------------------------------------------------------------
# Edit 3: "Since classes should be responsible for their own properties"
#         -> "Using the design principle that classes should be
#         responsible for their own properties" (bolding the
#         "classes should be responsible for their own properties" part)
# ---------------------------------------------------------------------
$r3 = Find-Unique "Since classes should be responsible for their own properties, it "
$sinceStart = $r3.Start
$sinceLen = "Since".Length
$classesStart = $sinceStart + $sinceLen + 1   # skip "Since "
$classesLen = "classes should be responsible for their own properties".Length

# Bold the "classes should be responsible..." phrase first (offsets are
# still valid because this step does not change the character count).
$d.Range($classesStart, $classesStart + $classesLen).Bold = 1

# Now swap out the leading "Since" for the new lead-in text.
$d.Range($sinceStart, $sinceStart + $sinceLen).Text = "Using the design principle that"

# ---------------------------------------------------------------------
# Edit 4: append the new FF / DRY / grouping-encapsulation sentence to
#         the end of the "Beating up the Zombies" paragraph.
# ---------------------------------------------------------------------
$r4 = Find-Unique "at all times."
$insertStart = $r4.End

$segments4 = @(
    @{ bold = $false; text = " Since the maximum number of legs or arms is 2 for each zombie, we use assertions to ensure that the value does not exceed 2, following the principle of " },
    @{ bold = $true;  text = "FF" },
    @{ bold = $false; text = " (Fail Fast). Since zombies can lose limbs to " },
    @{ bold = $false; text = "ANY" },
    @{ bold = $false; text = " attack that causes damage, to implement the losing limbs feature, we used the inherited method " },
    @{ bold = $false; text = "hurt(" },
    @{ bold = $false; text = ") from Actor class and override to add the feature" },
    @{ bold = $false; text = " using the principle " },
    @{ bold = $true;  text = "DRY" },
    @{ bold = $false; text = ". Keeping all the methods and attributes within the Zombie class, we followed the design principle of " },
    @{ bold = $true;  text = "gr" },
    @{ bold = $true;  text = "oup" },
    @{ bold = $true;  text = "ing" },
    @{ bold = $true;  text = " elements that must depend on each other together inside an" },
    @{ bold = $true;  text = " " },
    @{ bold = $true;  text = "encapsulation boundary" },
    @{ bold = $false; text = " of a class." }
)

$newText4 = ""
foreach ($seg in $segments4) { $newText4 += $seg.text }

$insertPoint = $d.Range($insertStart, $insertStart)
$insertPoint.InsertAfter($newText4)

$pos = $insertStart
foreach ($seg in $segments4) {
    $len = $seg.text.Length
    if ($seg.bold) {
        $d.Range($pos, $pos + $len).Bold = 1
    }
    $pos = $pos + $len
}

# ---------------------------------------------------------------------
# Edit 5: wrap the "Crafting class has 2 attributes ... Don't Repeat
#         Yourself. " sentence with a "_GoBack" bookmark.
# ---------------------------------------------------------------------
$bookmarkText = "Crafting class has 2 attributes, club and mace which are both WeaponItem objects because the attributes and methods that club and mace are required to use are similar to objects of Weapon item class, using the principle of Don" + [char]0x2019 + "t Repeat Yourself. "
$r5 = Find-Unique $bookmarkText
$d.Bookmarks.Add("_GoBack", $r5)
